$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 177; this shifts the existing rows 177-275 down to 178-276
# and updates the used-range dimension automatically.
$ws.Rows(177).Insert()

# Populate the newly inserted row 177 with the new daily record.
$ws.Range("A177").Value = 4
$ws.Range("B177").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C177").Value = "Los Lagos"
$ws.Range("D177").Value = 44518
$ws.Range("E177").Value = 10
$ws.Range("F177").Value = 100114001
$ws.Range("G177").Value = "Papa"
$ws.Range("H177").Value = "Pehuenche"
$ws.Range("I177").Value = "1a nueva(o)"
$ws.Range("J177").Value = 250
$ws.Range("K177").Value = 14000
$ws.Range("L177").Value = 14000
$ws.Range("M177").Value = 14000
$ws.Range("N177").Value = "$/saco 25 kilos"
$ws.Range("O177").Value = "Región de La Araucanía"
$ws.Range("P177").Value = 560
$ws.Range("Q177").Value = 25
$ws.Range("R177").Value = "Hortaliza"
